$d = $word.ActiveDocument

$replacements = @(
    @{old = "30÷6="; new = "29÷2="},
    @{old = "22÷4="; new = "16÷5="},
    @{old = "54÷4="; new = "38÷9="},
    @{old = "66÷2="; new = "11÷7="},
    @{old = "58÷6="; new = "86÷7="},
    @{old = "92÷8="; new = "53÷4="},
    @{old = "90÷7="; new = "74÷4="},
    @{old = "88÷9="; new = "74÷3="},
    @{old = "10÷4="; new = "18÷3="},
    @{old = "49÷9="; new = "28÷6="},
    @{old = "41÷2="; new = "26÷7="},
    @{old = "57÷4="; new = "44÷9="},
    @{old = "41÷8="; new = "28÷3="},
    @{old = "63÷4="; new = "46÷2="},
    @{old = "85÷2="; new = "80÷2="},
    @{old = "72÷8="; new = "83÷8="},
    @{old = "97÷5="; new = "94÷5="},
    @{old = "74÷7="; new = "26÷8="},
    @{old = "93÷2="; new = "64÷3="},
    @{old = "75÷9="; new = "93÷4="},
    @{old = "12÷4="; new = "31÷3="},
    @{old = "49÷4="; new = "97÷3="},
    @{old = "64÷8="; new = "54÷4="},
    @{old = "38÷2="; new = "41÷8="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
